# Crop the bottom portion off the screenshot picture on slide 10 so that
# only the top table/register screenshot remains visible, and shrink the
# shape's height to match the new (cropped) image proportions.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item("Google Shape;250;p22")

# Crop off the bottom ~42.366% of the source image (keeps only the top
# "edge node register" screenshot, drops the "application instance" one).
$shp.PictureFormat.CropBottom = 176.03

# The shape keeps the same width/position but its height shrinks to match
# the now-cropped picture (5458200 x 2045450 EMU == 429.7795 x 161.0591 pt).
$shp.Height = 161.0591
